$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.719.98'
$ws.Range("E2").Value = '  -0.32%  '

$ws.Range("D3").Value = '2.528.19'
$ws.Range("E3").Value = '  -1.56%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '309.22'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.18%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '100.53'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.29%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.568'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -1.15%  '

$ws.Range("E8").Value = '  +0.07%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.521'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.42%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '35.59'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.59%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0804'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.08%  '

$ws.Range("E12").Value = '  -1.89%  '

$ws.Range("E13").Value = '  +1.04%  '

$ws.Range("D14").Value = '2.915.27'
$ws.Range("E14").Value = '  -1.65%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '15.30'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -3.94%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '2.529.13'
$ws.Range("E16").Value = '  -3.15%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.811'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -4.22%  '

$ws.Range("D18").Value = '42.691.01'
$ws.Range("E18").Value = '  -0.45%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.68'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.31%  '

$ws.Range("D20").Value = '0.0₃0949'
$ws.Range("E20").Value = '  -1.40%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '12.26'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -2.33%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '69.29'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.39%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '242.56'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -2.93%  '

$ws.Range("E24").Value = '  -3.30%  '

$ws.Range("E25").Value = '  -2.70%  '

$ws.Range("E26").Value = '  +0.05%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '25.40'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -6.47%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.35'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -2.05%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '10.15'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.14%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '38.41'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -3.44%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '159.56'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.27%  '

$ws.Range("E32").Value = '  -0.86%  '

$ws.Range("E33").Value = '  +9.31%  '

$ws.Range("E34").Value = '  -1.28%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0782'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -2.48%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '18.47'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.08%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.12'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -7.46%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.96'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -7.26%  '

$ws.Range("E39").Value = '  -1.49%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.118'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.70%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '4.22'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +2.79%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '22.10'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -5.05%  '

$ws.Range("E43").Value = '  +0.13%  '

$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '3.30'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +1.80%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0300'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.88%  '

$ws.Range("D46").Value = '1.999.78'
$ws.Range("E46").Value = '  +0.01%  '

$ws.Range("E47").Value = '  -1.62%  '

$ws.Range("D48").Value = '2.769.44'
$ws.Range("E48").Value = '  -1.64%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.189'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -4.24%  '

$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '79.28'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -3.34%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '100.85'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.39%  '
